# Disable smart-quote autocorrect so straight quotes/apostrophes in the
# Zulu replacement text survive verbatim (defensive; we also avoid the
# Find/Replace auto-format path entirely by using Range.Text assignment).
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

$d = $word.ActiveDocument

# Paragraph 4
$p = $d.Paragraphs(4)
$r = $p.Range
[void]$r.Find.Execute('Welcome', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Siyakwamukela'

# Paragraph 6
$p = $d.Paragraphs(6)
$r = $p.Range
[void]$r.Find.Execute('Welcome to ParentText South Africa! ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Siyakwamukela kwi-ParentText South Africa! '

# Paragraph 8
$p = $d.Paragraphs(8)
$r = $p.Range
[void]$r.Find.Execute('ParentText is a chatbot service that helps you complete your parenting goals using a curriculum designed by Parenting for Lifelong Health with UNICEF and the World Health Organization, and tested all over the world. This programme works! ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'I-ParentText yinkundla yokuxoxa ekusiza ekutheni ufeze izinhloso zakho zokuba umzali ngokusebenzisa ikharikhulamu eyakhiwe yi-Parenting for Lifelong Health ibambisene no-UNICEF kanye ne-World Health Organisation, futhi ehlolwe emhlabeni wonke. Loluhlelo luyasebenza! '

# Paragraph 10
$p = $d.Paragraphs(10)
$r = $p.Range
[void]$r.Find.Execute('Being here shows how much you care about providing the best support for your teen. Halala!', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukuba khona kwakho lana kutshengisa ukuthi ukukhathalele ukunikeza umntwana wakho ukwesekwa okuvelele. Halala!'

# Paragraph 12
$p = $d.Paragraphs(12)
$r = $p.Range
[void]$r.Find.Execute('Remember: it is what you do with your teen that will make a difference. ParentText will provide you with tips and skills to help you with your relationship with your teen, but it is up to you to put these tips into practice!', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Khumbula: ilokho okwenza nomntwana wakho okwenza umehluko. I-ParentText izokuhlinzeka ngamacebo namakhono azokusiza ebudlelwaneni bakho nomntwana wakho, kodwa kukuwe ukuwasebenzisa lamacebo!'

# Paragraph 14
$p = $d.Paragraphs(14)
$r = $p.Range
[void]$r.Find.Execute('I’m Ayanda, your guide. I may look like a human, but I’m actually a robot produced by Parenting for Lifelong Health and UNICEF to help you learn. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ngingu Ayanda, umhlahlandlela wakho. Ngingabukeka ngathi ngingumuntu, kwodwa ngiyirobhothi elizokusiza eLakhiwe yi-Parenting for Lifelong Health no UNICEF. '

# Paragraph 16
$p = $d.Paragraphs(16)
$r = $p.Range
[void]$r.Find.Execute('Today, I’m going to explain how to use ParentText. Together we will review: ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Namhlanje ngizokuchazela ukuthi isetshenziswa kanjani i-ParentText. Ndawonye sizobuyekeza: '

# Paragraph 17
$p = $d.Paragraphs(17)
$r = $p.Range
[void]$r.Find.Execute('How to earn your Positive Parenting Trophy', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukuthi uyizuza kanjani iNdondo yokuba uMzali oMuhle'

# Paragraph 18
$p = $d.Paragraphs(18)
$r = $p.Range
[void]$r.Find.Execute('How to make progress in your parenting goals', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukuthi uzithuthukisa kanjani izinhloso zakho zobuzali'

# Paragraph 19
$p = $d.Paragraphs(19)
$r = $p.Range
[void]$r.Find.Execute('How to track your progress', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukuthi uyilandelela kanjani inqubekela phambili yakho'

# Paragraph 20
$p = $d.Paragraphs(20)
$r = $p.Range
[void]$r.Find.Execute('How to get help with this course', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukuthi ulithola kanjani usizo kulesisifundo'

# Paragraph 21
$p = $d.Paragraphs(21)
$r = $p.Range
[void]$r.Find.Execute('Accessing support to troubleshoot common parenting challenges, and', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukuthola ukwesekwa ekuxazululeni izinselelo zobuzali ezivamile, ne'

# Paragraph 22
$p = $d.Paragraphs(22)
$r = $p.Range
[void]$r.Find.Execute('Resources available to you in an emergency or crisis. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Izinsiza ongazithola lapho ubhekana nesimo esiphuthumayo noma inhlekele. '

# Paragraph 25
$p = $d.Paragraphs(25)
$r = $p.Range
[void]$r.Find.Execute('Earning Your Positive Parenting Trophy', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukuthola iNdondo Yakho yokuba uMzali Omuhle'

# Paragraph 27
$p = $d.Paragraphs(27)
$r = $p.Range
[void]$r.Find.Execute('First, How to complete the course and earn your Positive Parenting Trophy.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Okokuqala, Usiqeda kanjani isifundo bese uthola iNdondo Yakho yokuba uMzali Omuhle.'
$r = $p.Range
[void]$r.Find.Execute(' Depending on your profile, this course is divided into 8 or 9 parenting goals to help you improve your relationship with your teen and help them thrive. Each goal takes between two to four days to complete, with each day bringing a new lesson.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Ngokuncika kwiphrofiyili yakho, lesisifundo sicazwe ngezinghloso zobuzali eziwu 8 noma 9 ezizokusiza ekuthuthukiseni ubudlelwane bakho nomntwana wakho futhi kusize ekutheni umntwana wakho aqhakaze. Inhloso ngayinye ithatha phakathi kwezinsuku ezimbili ukuya kwezine ukuthi uyiqede, usuku nosuku luza nesifundo esisha.'
$r = $p.Range
[void]$r.Find.Execute(' You can choose which goals you want to work on first, but each goal must be completed in order to complete the course and earn your Positive Parenting Trophy.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Ungakhetha ukuthi iziphi izinjongo ofuna ukuqala ngazo, kodwa injongo nenjongo kumele iqedwe khona kuzoqedeka isifundo bese uthola iNdondo Yakho yokuba uMzali Omuhle.'

# Paragraph 30
$p = $d.Paragraphs(30)
$r = $p.Range
[void]$r.Find.Execute('Completing a Goal', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukuqeda iNjongo'

# Paragraph 31
$p = $d.Paragraphs(31)
$r = $p.Range
[void]$r.Find.Execute(' Now, let''s learn how to complete each of the goals in the programme.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Manje masifunde ukuthi siziqeda kanjani lezinjongo ezikuloluhlelo.'
$r = $p.Range
[void]$r.Find.Execute(' Soon, you''ll be asked to select your first parenting goal. Once a goal has been selected, you will complete lessons to gain new skills.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Maduze, uzocelwa ukuthi ukhethe inhloso yakho yokuqala. Uma inhloso isikhethiwe, uzoqedela izifundondo khona uzozuza amakhono amasha.'
$r = $p.Range
[void]$r.Find.Execute(' A new skill is offered every day. Most skills take less than 5 minutes to complete. If you don''t manage to complete a skill assigned, I will ask if you want to complete it the next day.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Ikhono elisha litholakala nsukuzonke. Kuthatha ngaphansi kwemizuzu emihlanu ukuqeda iningi lalamakhono. Uma ungakwazanga ukuqeda ikhono owabelwe lona, ngizokubuza ukuthi uyafuna yini ukuliqedela ngakusasa.'
$r = $p.Range
[void]$r.Find.Execute(' You must complete all the skills within a goal to earn a badge. After you have completed one parenting goal, you can select another. Once all of your goal badges are earned, you will have completed the course, and will receive the Positive Parenting Trophy.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Kumele uwaqede wonke amakhono akuleyonhloso yesifundo ukuze uthole ibheji. Uma usuyiqedile inhloso eyodwa yobuzali, ungakhetha enye. Uma usuwathole wonke amabheji enhloso, uzobe ususiqedile isifundo bese uthola iNdondo yakho yoBuzali Obuhle.'

# Paragraph 34
$p = $d.Paragraphs(34)
$r = $p.Range
[void]$r.Find.Execute('Tracking Progress', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ukulandelela inqubekelaphambili'

# Paragraph 35
$p = $d.Paragraphs(35)
$r = $p.Range
[void]$r.Find.Execute('As you move through each day''s lesson, you''ll receive updates on your progress that look like this: . These check marks tell you how far along you are on the day''s lesson.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Uma uqhubeka nesifundo sosuku, uzothola izibuyekezo ngenqubekelaphambili yakho ebukeka kanje:  Lezizimpawu zikutshela ukuthi usuhambe kanganani esifundweni sakho sosuku.'
$r = $p.Range
[void]$r.Find.Execute(' If you want to see how far you have progressed with your parenting goals, you can check on your goal progress through the Main Menu. To access the menu, type "Menu" at any time.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Uma ufuna ukubona ukuthi usuhambe kangakanani ngezihloso zakho zobuzali, ungabheka kwi- Main Menu. Ukuthola imenyu, bhala u ''Menu" noma yinini.'
$r = $p.Range
[void]$r.Find.Execute(' The first menu option is labelled "track my progress."', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Inketho yokuqala kwimenyu imakwe ngokuthi "bheka inqubekelaphambili yami".'
$r = $p.Range
[void]$r.Find.Execute(' Here you can see your progress, review the goals you have achieved, and those that are still incomplete.', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = ' Kulapha ke lapho ubona inqubekela phambili yakho, ubuyekeze izinhloso ozifezile nalezo ongakazifezi.'

# Paragraph 38
$p = $d.Paragraphs(38)
$r = $p.Range
[void]$r.Find.Execute('Menu ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Menyu '

# Paragraph 39
$p = $d.Paragraphs(39)
$r = $p.Range
[void]$r.Find.Execute('The menu contains other features that might help you, too. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Imenyu iqukethe ezinye izici ezingakusiza. '

# Paragraph 40
$p = $d.Paragraphs(40)
$r = $p.Range
[void]$r.Find.Execute('In addition to tracking your progress, you can also: ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ngokwengeziwe ekubhekeni inqubekelaphambili yakho, unga: '

# Paragraph 41
$p = $d.Paragraphs(41)
$r = $p.Range
[void]$r.Find.Execute('Share ParentText with a friend and help them enroll. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Yabelana nomngani i-ParentText bese uyamsiza abhalise. '

# Paragraph 42
$p = $d.Paragraphs(42)
$r = $p.Range
[void]$r.Find.Execute('Change your settings, like how you receive the messages, when you receive notifications, or adjust information about yourself and your teen to get the best support. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Shintsha ama-settings wakho, njengendlela yokuthi uyithola kanjani imiyalezo, uzithola nini izaziso, noma ukulungisa imininingwane engawe neyomntwana wakho khona nizothola usizo. '

# Paragraph 43
$p = $d.Paragraphs(43)
$r = $p.Range
[void]$r.Find.Execute('Access a list of activities that you can complete with your teen to build your relationship. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Thola uhla lwezinto zokwenza ongazenza nomtwana wakho ukwakha ubudlelwane benu. '

# Paragraph 44
$p = $d.Paragraphs(44)
$r = $p.Range
[void]$r.Find.Execute('Review this onboarding guide and receive support navigating ParentText. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Buyekeza lomhlahlandlela wokuqondisa nokuthola usizo lokusebenzisa i-ParentText. '

# Paragraph 45
$p = $d.Paragraphs(45)
$r = $p.Range
[void]$r.Find.Execute('And get help troubleshooting difficult challenges with your teen. Let''s learn more about this feature now. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Bese uthola usizo ekuxazululeni izinselelo ezinzima nomntwana wakho. Manje ake sifunde kabanzi ngalesi sici. '

# Paragraph 48
$p = $d.Paragraphs(48)
$r = $p.Range
[void]$r.Find.Execute('Troubleshooting', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Xazulula izinkinga'

# Paragraph 49
$p = $d.Paragraphs(49)
$r = $p.Range
[void]$r.Find.Execute('Parenting can be difficult. Though challenges feel unique to you, they are often more common than you think. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Ubuzali bungaba nzima. Izinselelo ezinzima kakhulu zizwakala ngathi uwe wedwa obhekene naza, kanti zivamile kunalokho okucabangayo. '

# Paragraph 50
$p = $d.Paragraphs(50)
$r = $p.Range
[void]$r.Find.Execute('As you begin to achieve goals in this programme, I will check in on how things are going with your teen. If they aren’t going well, I might offer support. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Lapho uqala ukufeza izinhloso zakho kuloluhlelo, ngizobheka ukuthi izinto zihamba kanjani phakathi kwakho nomntwana wakho. Uma zingahambi kahle, kungenzeka ngikusize. '

# Paragraph 51
$p = $d.Paragraphs(51)
$r = $p.Range
[void]$r.Find.Execute('When you share with me your challenges, I will offer practical solutions to help you succeed. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Uma ungitshela izinkinga zakho, ngizokunikeza izixazululo ezingokoqobo ezizokusiza ukuthi uphumelele. '

# Paragraph 52
$p = $d.Paragraphs(52)
$r = $p.Range
[void]$r.Find.Execute('You don’t have to wait on me to offer support, though. You can also access troubleshooting support through the Main Menu at any time. ', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$r.Text = 'Kodwa, awidingi ukulinda kuze kube yimina oza kuwe nosizo. Nawe ungaluthola usizo lokuxazululaizinkinga lapha kwiMenyu nanoma yinini. '
